# Insert a new data row right before the current row 128 (shifts existing
# rows 128..189 down to 129..190, matching the diff's dimension change from
# A1:R189 to A1:R190 and the observed "each row's data moved to row+1").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(128).Insert()

# Populate the newly inserted row 128 with the new record.
$ws.Cells.Item(128, 1).Value  = 8
$ws.Cells.Item(128, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(128, 3).Value  = "Coquimbo"
$ws.Cells.Item(128, 4).Value  = 45089
$ws.Cells.Item(128, 5).Value  = 4
$ws.Cells.Item(128, 6).Value  = 100112052
$ws.Cells.Item(128, 7).Value  = "Albahaca"
$ws.Cells.Item(128, 8).Value  = "Sin especificar"
$ws.Cells.Item(128, 9).Value  = "Primera"
$ws.Cells.Item(128, 10).Value = 800
$ws.Cells.Item(128, 11).Value = 3000
$ws.Cells.Item(128, 12).Value = 3500
$ws.Cells.Item(128, 13).Value = 3250
$ws.Cells.Item(128, 14).Value = "$/paquete"
$ws.Cells.Item(128, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(128, 16).Value = 3250
$ws.Cells.Item(128, 17).Value = 1
$ws.Cells.Item(128, 18).Value = "Hortaliza"
